# Update TPM-derived values in the LR-pairs sheet (Hbegf-Erbb2)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns G..T correspond to the numeric ligand/receptor/edge expression
# statistics that were recomputed with the new TPM values.
# Row => ordered values for columns G,H,I,J,K,L,M,N,O,P,Q,R,S,T
$data = @{
    2  = @(7.723979, 23.171937, 0.471042132528101, 0.471042132528101, 3, 1, 2.451416666666667, 7.35425, 0.2191928499183569, 0.2191928499183569, 18.93469085358333, 170.41221768225, 0.1032490674604548, 0.1032490674604548)
    3  = @(7.723979, 23.171937, 0.471042132528101, 0.471042132528101, 3, 1, 4.973328333333334, 14.919985, 0.4446889938320204, 0.4446889938320204, 38.41388360677167, 345.724952460945, 0.2094672519664104, 0.2094672519664104)
    4  = @(7.723979, 23.171937, 0.471042132528101, 0.471042132528101, 3, 1, 3.75909, 11.27727, 0.3361181562496228, 0.3361181562496228, 29.03513221911, 261.31618997199, 0.1583258131012358, 0.1583258131012358)
    5  = @(4.034036666666666, 12.10211, 0.2460132574367717, 0.2460132574367717, 3, 1, 2.451416666666667, 7.35425, 0.2191928499183569, 0.2191928499183569, 9.88910471861111, 89.0019424675, 0.0539243470152644, 0.0539243470152644)
    6  = @(4.034036666666666, 12.10211, 0.2460132574367717, 0.2460132574367717, 3, 1, 4.973328333333334, 14.919985, 0.4446889938320204, 0.4446889938320204, 20.06258885203889, 180.56329966835, 0.1093993879188958, 0.1093993879188958)
    7  = @(4.034036666666666, 12.10211, 0.2460132574367717, 0.2460132574367717, 3, 1, 3.75909, 11.27727, 0.3361181562496228, 0.3361181562496228, 15.1643068933, 136.4787620397, 0.08268952250261151, 0.08268952250261151)
    8  = @(4.639623666666666, 13.918871, 0.2829446100351274, 0.2829446100351274, 3, 1, 2.451416666666667, 7.35425, 0.2191928499183569, 0.2191928499183569, 11.37365078352778, 102.36285705175, 0.0620194354426377, 0.0620194354426377)
    9  = @(4.639623666666666, 13.918871, 0.2829446100351274, 0.2829446100351274, 3, 1, 4.973328333333334, 14.919985, 0.4446889938320204, 0.4446889938320204, 23.07437183743722, 207.669346536935, 0.1258223539467142, 0.1258223539467142)
    10 = @(4.639623666666666, 13.918871, 0.2829446100351274, 0.2829446100351274, 3, 1, 3.75909, 11.27727, 0.3361181562496228, 0.3361181562496228, 17.44076292913, 156.96686636217, 0.09510282064577555, 0.09510282064577555)
}

$columns = @("G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $ws.Range("$($columns[$i])$row").Value = $values[$i]
    }
}
